$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width ---
$ws.Columns("B").ColumnWidth = 45.85546875

# --- Row 8 ---
$ws.Range("A8").Value = 45668
$ws.Range("A8").NumberFormat = "m/d/yyyy"
$ws.Range("A8").ShrinkToFit = $false
$ws.Range("B8").ShrinkToFit = $false

# --- Row 9 ---
$ws.Range("A9").Value = 45669
$ws.Range("A9").NumberFormat = "m/d/yyyy"
$ws.Range("A9").ShrinkToFit = $false
$ws.Range("B9").ShrinkToFit = $false

# --- Row 10 ---
$ws.Range("A10").Value = 45670
$ws.Range("A10").NumberFormat = "m/d/yyyy"
$ws.Range("A10").ShrinkToFit = $false
$ws.Range("B10").ShrinkToFit = $false
$ws.Range("B10").Value = "pongal holidays"

# --- Row 11 ---
$ws.Range("A11").Value = 45671
$ws.Range("A11").NumberFormat = "m/d/yyyy"
$ws.Range("A11").ShrinkToFit = $false
$ws.Range("B11").ShrinkToFit = $false

# --- Row 12 ---
$ws.Range("A12").Value = 45672
$ws.Range("A12").NumberFormat = "m/d/yyyy"
$ws.Range("A12").ShrinkToFit = $false
$ws.Range("B12").ShrinkToFit = $false

# --- Row 13 ---
$ws.Range("A13").Value = 45673
$ws.Range("A13").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B13").Value = "project work,dataiku instalation,dataiku video2"
$ws.Range("D13").Value = "completed"

# --- Selection ---
$ws.Range("B10").Select()
